$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Rename every "TOTAL" subtotal label in column B to "Total" ---
$totalRows = @(4,7,10,13,16,19,22,25,28,31,34,37,40,43,46,49,52,55,58,61,64,66,69,72,75,78,81,84,87,90)
foreach ($r in $totalRows) {
    $ws.Cells.Item($r, 2).Value = "Total"
}

# --- Grand total row 91, column A ---
$ws.Range("A91").Value = "Total"

# --- Title-case a handful of state names in column A ---
$ws.Range("A17").Value = "Ciudad De México"
$ws.Range("A20").Value = "Coahuila De Zaragoza"
$ws.Range("A29").Value = "Estado De México"
$ws.Range("A44").Value = "Michoacán De Ocampo"
$ws.Range("A85").Value = "Veracruz De Ignacio De La Llave"

# --- Remove trailing footnote rows (93-97) ---
$ws.Range("A93:A97").EntireRow.Delete()

# --- Fix the sheet dimension to reflect the new used range ---
$ws.Range("A1:D91").Select()
